$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.0134664339880133
$ws.Range("J2").Value = 0.0134664339880133
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.224187333333333
$ws.Range("N2").Value = 3.672562
$ws.Range("O2").Value = 0.2801938168429309
$ws.Range("P2").Value = 0.2801938168429309
$ws.Range("Q2").Value = 0.1021351734073333
$ws.Range("R2").Value = 0.919216560666
$ws.Range("S2").Value = 0.003773211538364817
$ws.Range("T2").Value = 0.003773211538364818

# Row 3
$ws.Range("I3").Value = 0.0134664339880133
$ws.Range("J3").Value = 0.0134664339880133
$ws.Range("N3").Value = 4.922466999999999
$ws.Range("O3").Value = 0.3755538550508803
$ws.Range("P3").Value = 0.3755538550508803
$ws.Range("S3").Value = 0.005057371197986595
$ws.Range("T3").Value = 0.005057371197986595

# Row 4
$ws.Range("I4").Value = 0.0134664339880133
$ws.Range("J4").Value = 0.0134664339880133
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.254101
$ws.Range("N4").Value = 0.762303
$ws.Range("O4").Value = 0.05815901464994103
$ws.Range("P4").Value = 0.05815901464994102
$ws.Range("Q4").Value = 0.02119990053099999
$ws.Range("R4").Value = 0.190799104779
$ws.Range("S4").Value = 0.0007831945315913292
$ws.Range("T4").Value = 0.0007831945315913292

# Row 5
$ws.Range("I5").Value = 0.0134664339880133
$ws.Range("J5").Value = 0.0134664339880133
$ws.Range("M5").Value = 0.8417103333333333
$ws.Range("N5").Value = 2.525131
$ws.Range("O5").Value = 0.192651912457409
$ws.Range("P5").Value = 0.192651912457409
$ws.Range("Q5").Value = 0.07022473482033333
$ws.Range("R5").Value = 0.6320226133829999
$ws.Range("S5").Value = 0.002594334261772215
$ws.Range("T5").Value = 0.002594334261772215

# Row 6
$ws.Range("I6").Value = 0.0134664339880133
$ws.Range("J6").Value = 0.0134664339880133
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4082523333333333
$ws.Range("N6").Value = 1.224757
$ws.Range("O6").Value = 0.09344140099883881
$ws.Range("P6").Value = 0.0934414009988388
$ws.Range("Q6").Value = 0.03406090042233333
$ws.Range("R6").Value = 0.3065481038009999
$ws.Range("S6").Value = 0.001258322458298343
$ws.Range("T6").Value = 0.001258322458298343

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.699187666666667
$ws.Range("H7").Value = 14.097563
$ws.Range("I7").Value = 0.758486659760196
$ws.Range("J7").Value = 0.758486659760196
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.224187333333333
$ws.Range("N7").Value = 3.672562
$ws.Range("O7").Value = 0.2801938168429309
$ws.Range("P7").Value = 0.2801938168429309
$ws.Range("Q7").Value = 5.752686018489555
$ws.Range("R7").Value = 51.77417416640601
$ws.Range("S7").Value = 0.2125232722226548
$ws.Range("T7").Value = 0.2125232722226548

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.699187666666667
$ws.Range("H8").Value = 14.097563
$ws.Range("I8").Value = 0.758486659760196
$ws.Range("J8").Value = 0.758486659760196
$ws.Range("N8").Value = 4.922466999999999
$ws.Range("O8").Value = 0.3755538550508803
$ws.Range("P8").Value = 0.3755538550508803
$ws.Range("Q8").Value = 7.710532071991222
$ws.Range("R8").Value = 69.39478864792099
$ws.Range("S8").Value = 0.2848525890776071
$ws.Range("T8").Value = 0.284852589077607

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.699187666666667
$ws.Range("H9").Value = 14.097563
$ws.Range("I9").Value = 0.758486659760196
$ws.Range("J9").Value = 0.758486659760196
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.254101
$ws.Range("N9").Value = 0.762303
$ws.Range("O9").Value = 0.05815901464994103
$ws.Range("P9").Value = 0.05815901464994102
$ws.Range("Q9").Value = 1.194068285287667
$ws.Range("R9").Value = 10.746614567589
$ws.Range("S9").Value = 0.04411283675677807
$ws.Range("T9").Value = 0.04411283675677807

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.699187666666667
$ws.Range("H10").Value = 14.097563
$ws.Range("I10").Value = 0.758486659760196
$ws.Range("J10").Value = 0.758486659760196
$ws.Range("M10").Value = 0.8417103333333333
$ws.Range("N10").Value = 2.525131
$ws.Range("O10").Value = 0.192651912457409
$ws.Range("P10").Value = 0.192651912457409
$ws.Range("Q10").Value = 3.955354817305889
$ws.Range("R10").Value = 35.598193355753
$ws.Range("S10").Value = 0.1461239055762339
$ws.Range("T10").Value = 0.1461239055762338

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.699187666666667
$ws.Range("H11").Value = 14.097563
$ws.Range("I11").Value = 0.758486659760196
$ws.Range("J11").Value = 0.758486659760196
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.4082523333333333
$ws.Range("N11").Value = 1.224757
$ws.Range("O11").Value = 0.09344140099883881
$ws.Range("P11").Value = 0.0934414009988388
$ws.Range("Q11").Value = 1.918454329687889
$ws.Range("R11").Value = 17.266088967191
$ws.Range("S11").Value = 0.0708740561269223
$ws.Range("T11").Value = 0.07087405612692228

# Row 12
$ws.Range("G12").Value = 0.1492686666666667
$ws.Range("H12").Value = 0.447806
$ws.Range("I12").Value = 0.02409316256721636
$ws.Range("J12").Value = 0.02409316256721636
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.224187333333333
$ws.Range("N12").Value = 3.672562
$ws.Range("O12").Value = 0.2801938168429309
$ws.Range("P12").Value = 0.2801938168429309
$ws.Range("Q12").Value = 0.1827328109968889
$ws.Range("R12").Value = 1.644595298972
$ws.Range("S12").Value = 0.00675075517952558
$ws.Range("T12").Value = 0.006750755179525578

# Row 13
$ws.Range("G13").Value = 0.1492686666666667
$ws.Range("H13").Value = 0.447806
$ws.Range("I13").Value = 0.02409316256721636
$ws.Range("J13").Value = 0.02409316256721636
$ws.Range("N13").Value = 4.922466999999999
$ws.Range("O13").Value = 0.3755538550508803
$ws.Range("P13").Value = 0.3755538550508803
$ws.Range("Q13").Value = 0.2449233619335556
$ws.Range("R13").Value = 2.204310257402
$ws.Range("S13").Value = 0.00904828008248567
$ws.Range("T13").Value = 0.009048280082485666

# Row 14
$ws.Range("G14").Value = 0.1492686666666667
$ws.Range("H14").Value = 0.447806
$ws.Range("I14").Value = 0.02409316256721636
$ws.Range("J14").Value = 0.02409316256721636
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.254101
$ws.Range("N14").Value = 0.762303
$ws.Range("O14").Value = 0.05815901464994103
$ws.Range("P14").Value = 0.05815901464994102
$ws.Range("Q14").Value = 0.03792931746866667
$ws.Range("R14").Value = 0.341363857218
$ws.Range("S14").Value = 0.001401234594710147
$ws.Range("T14").Value = 0.001401234594710147

# Row 15
$ws.Range("G15").Value = 0.1492686666666667
$ws.Range("H15").Value = 0.447806
$ws.Range("I15").Value = 0.02409316256721636
$ws.Range("J15").Value = 0.02409316256721636
$ws.Range("M15").Value = 0.8417103333333333
$ws.Range("N15").Value = 2.525131
$ws.Range("O15").Value = 0.192651912457409
$ws.Range("P15").Value = 0.192651912457409
$ws.Range("Q15").Value = 0.1256409791762222
$ws.Range("R15").Value = 1.130768812586
$ws.Range("S15").Value = 0.004641593845721491
$ws.Range("T15").Value = 0.004641593845721489

# Row 16
$ws.Range("G16").Value = 0.1492686666666667
$ws.Range("H16").Value = 0.447806
$ws.Range("I16").Value = 0.02409316256721636
$ws.Range("J16").Value = 0.02409316256721636
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.4082523333333333
$ws.Range("N16").Value = 1.224757
$ws.Range("O16").Value = 0.09344140099883881
$ws.Range("P16").Value = 0.0934414009988388
$ws.Range("Q16").Value = 0.06093928146022222
$ws.Range("R16").Value = 0.548453533142
$ws.Range("S16").Value = 0.002251298864773477
$ws.Range("T16").Value = 0.002251298864773476

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.263591
$ws.Range("H17").Value = 3.790773
$ws.Range("I17").Value = 0.2039537436845743
$ws.Range("J17").Value = 0.2039537436845743
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.224187333333333
$ws.Range("N17").Value = 3.672562
$ws.Range("O17").Value = 0.2801938168429309
$ws.Range("P17").Value = 0.2801938168429309
$ws.Range("Q17").Value = 1.546872096714
$ws.Range("R17").Value = 13.921848870426
$ws.Range("S17").Value = 0.05714657790238567
$ws.Range("T17").Value = 0.05714657790238567

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1.263591
$ws.Range("H18").Value = 3.790773
$ws.Range("I18").Value = 0.2039537436845743
$ws.Range("J18").Value = 0.2039537436845743
$ws.Range("N18").Value = 4.922466999999999
$ws.Range("O18").Value = 0.3755538550508803
$ws.Range("P18").Value = 0.3755538550508803
$ws.Range("Q18").Value = 2.073328332999
$ws.Range("R18").Value = 18.659954996991
$ws.Range("S18").Value = 0.076595614692801
$ws.Range("T18").Value = 0.07659561469280099

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1.263591
$ws.Range("H19").Value = 3.790773
$ws.Range("I19").Value = 0.2039537436845743
$ws.Range("J19").Value = 0.2039537436845743
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.254101
$ws.Range("N19").Value = 0.762303
$ws.Range("O19").Value = 0.05815901464994103
$ws.Range("P19").Value = 0.05815901464994102
$ws.Range("Q19").Value = 0.321079736691
$ws.Range("R19").Value = 2.889717630219
$ws.Range("S19").Value = 0.01186174876686147
$ws.Range("T19").Value = 0.01186174876686147

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1.263591
$ws.Range("H20").Value = 3.790773
$ws.Range("I20").Value = 0.2039537436845743
$ws.Range("J20").Value = 0.2039537436845743
$ws.Range("M20").Value = 0.8417103333333333
$ws.Range("N20").Value = 2.525131
$ws.Range("O20").Value = 0.192651912457409
$ws.Range("P20").Value = 0.192651912457409
$ws.Range("Q20").Value = 1.063577601807
$ws.Range("R20").Value = 9.572198416263001
$ws.Range("S20").Value = 0.03929207877368143
$ws.Range("T20").Value = 0.03929207877368143

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1.263591
$ws.Range("H21").Value = 3.790773
$ws.Range("I21").Value = 0.2039537436845743
$ws.Range("J21").Value = 0.2039537436845743
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.4082523333333333
$ws.Range("N21").Value = 1.224757
$ws.Range("O21").Value = 0.09344140099883881
$ws.Range("P21").Value = 0.0934414009988388
$ws.Range("Q21").Value = 0.515863974129
$ws.Range("R21").Value = 4.642775767161
$ws.Range("S21").Value = 0.01905772354884469
$ws.Range("T21").Value = 0.01905772354884469
